$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.2062081610468891
$ws.Range("C2").Value = 0.4981741733236962
$ws.Range("D2").Value = 0.4649261630253289
$ws.Range("E2").Value = 0.6818549428033274
$ws.Range("F2").Value = 0.657108150032778
$ws.Range("G2").Value = 46

$ws.Range("B3").Value = -0.09587357185906692
$ws.Range("C3").Value = 0.5000820061429816
$ws.Range("D3").Value = 0.4618601616711875
$ws.Range("E3").Value = 0.6796029441307531
$ws.Range("F3").Value = 0.6804089494737995
$ws.Range("G3").Value = 45

$ws.Range("B4").Value = -0.1372403718692351
$ws.Range("C4").Value = 0.4340014692588762
$ws.Range("D4").Value = 0.3663092181924378
$ws.Range("E4").Value = 0.6052348454876321
$ws.Range("F4").Value = 0.5962844087863356
$ws.Range("G4").Value = 44

$ws.Range("B5").Value = -0.08352607004722953
$ws.Range("C5").Value = 0.471891385230405
$ws.Range("D5").Value = 0.4001742410426709
$ws.Range("E5").Value = 0.6325932666750974
$ws.Range("F5").Value = 0.6344757561618594
$ws.Range("G5").Value = 43

$ws.Range("B6").Value = -0.1513393225090566
$ws.Range("C6").Value = 0.4428429211668459
$ws.Range("D6").Value = 0.3795138264528281
$ws.Range("E6").Value = 0.616046935267783
$ws.Range("F6").Value = 0.6044071860487836
$ws.Range("G6").Value = 42

$ws.Range("B7").Value = -0.1148246166180287
$ws.Range("C7").Value = 0.4704721399279175
$ws.Range("D7").Value = 0.4088047900674556
$ws.Range("E7").Value = 0.6393784404149514
$ws.Range("F7").Value = 0.6367971418930269
$ws.Range("G7").Value = 41

$ws.Range("B8").Value = -0.1576462645105448
$ws.Range("C8").Value = 0.4410790206399262
$ws.Range("D8").Value = 0.3843178504030022
$ws.Range("E8").Value = 0.6199337467850917
$ws.Range("F8").Value = 0.6071923665011824
$ws.Range("G8").Value = 40

$ws.Range("B9").Value = -0.1265473388711615
$ws.Range("C9").Value = 0.4734734975690142
$ws.Range("D9").Value = 0.415038962515679
$ws.Range("E9").Value = 0.6442351764035235
$ws.Range("F9").Value = 0.6399417039254013
$ws.Range("G9").Value = 39

$ws.Range("B10").Value = -0.181494460453781
$ws.Range("C10").Value = 0.4408424176267598
$ws.Range("D10").Value = 0.3840813816712615
$ws.Range("E10").Value = 0.6197429964681017
$ws.Range("F10").Value = 0.6005259724977672
$ws.Range("G10").Value = 38

$ws.Range("B11").Value = -0.1278002099571802
$ws.Range("C11").Value = 0.4693678920889392
$ws.Range("D11").Value = 0.4189108628417125
$ws.Range("E11").Value = 0.6472332368178511
$ws.Range("F11").Value = 0.6432423264545254
$ws.Range("G11").Value = 37

